# "print Subworks in Excel"
# Update the itemized-work list on the active sheet: correct the first two
# activity descriptions and append two new sub-work rows right below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing activity text (rows 19-20 of the "Itemizado trabajos" table)
$ws.Range("C19").Value = "  8 Cambio de Caps (60cu) "
$ws.Range("C20").Value = "Nueva67ffs8"

# Add the two new sub-work items under the existing ones
$ws.Range("C21").Value = "Nueva90"
$ws.Range("C22").Value = "Nueva7"

# Match the right-aligned style already used by the row above (s="21")
$ws.Range("C21").HorizontalAlignment = -4152
$ws.Range("C22").HorizontalAlignment = -4152

# Move the active selection to the last entry just filled in
$ws.Range("C22").Select()
